$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "57.225.94"
$ws.Cells.Item(2, 5).Value = "  +0.17%  "

$ws.Cells.Item(3, 4).Value = "2.425.99"
$ws.Cells.Item(3, 5).Value = "  -1.71%  "

$ws.Cells.Item(4, 5).Value = "  +0.20%  "

$c = $ws.Cells.Item(5, 4)
$c.Value = "'490.11"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.38%  "

$c = $ws.Cells.Item(6, 4)
$c.Value = "'156.08"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +1.44%  "

$ws.Cells.Item(7, 5).Value = "  +0.02%  "

$c = $ws.Cells.Item(8, 4)
$c.Value = "'0.610"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +19.03%  "

$ws.Cells.Item(9, 4).Value = "2.449.40"
$ws.Cells.Item(9, 5).Value = "  -1.00%  "

$c = $ws.Cells.Item(10, 4)
$c.Value = "'6.27"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +10.07%  "

$ws.Cells.Item(11, 5).Value = "  -0.59%  "

$ws.Cells.Item(12, 5).Value = "  -0.97%  "

$ws.Cells.Item(13, 5).Value = "  +0.83%  "

$ws.Cells.Item(14, 4).Value = "2.853.17"
$ws.Cells.Item(14, 5).Value = "  -1.50%  "

$ws.Cells.Item(15, 4).Value = "57.297.80"
$ws.Cells.Item(15, 5).Value = "  +0.15%  "

$c = $ws.Cells.Item(16, 4)
$c.Value = "'20.73"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -1.61%  "

$ws.Cells.Item(17, 5).Value = "  -2.99%  "

$ws.Cells.Item(18, 4).Value = "2.452.75"
$ws.Cells.Item(18, 5).Value = "  -1.02%  "

$ws.Cells.Item(19, 5).Value = "  +1.71%  "

$c = $ws.Cells.Item(20, 4)
$c.Value = "'325.08"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.47%  "

$ws.Cells.Item(21, 5).Value = "  -0.73%  "

$ws.Cells.Item(22, 5).Value = "  -0.18%  "

$c = $ws.Cells.Item(23, 4)
$c.Value = "'5.91"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +1.30%  "

$c = $ws.Cells.Item(24, 4)
$c.Value = "'58.10"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.26%  "

$c = $ws.Cells.Item(25, 4)
$c.Value = "'0.403"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -1.17%  "

$c = $ws.Cells.Item(26, 4)
$c.Value = "'0.998"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.33%  "

$c = $ws.Cells.Item(27, 4)
$c.Value = "'0.160"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -2.66%  "

$ws.Cells.Item(28, 4).Value = "2.551.39"
$ws.Cells.Item(28, 5).Value = "  -0.30%  "

$c = $ws.Cells.Item(29, 4)
$c.Value = "'7.30"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -3.35%  "

$ws.Cells.Item(30, 4).Value = "0.0₃0801"
$ws.Cells.Item(30, 5).Value = "  -2.36%  "

$ws.Cells.Item(31, 5).Value = "  +0.00%  "

$c = $ws.Cells.Item(32, 4)
$c.Value = "'150.80"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +0.22%  "

$c = $ws.Cells.Item(33, 4)
$c.Value = "'18.77"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +2.98%  "

$ws.Cells.Item(34, 5).Value = "  -0.11%  "

$ws.Cells.Item(35, 5).Value = "  +1.67%  "

$ws.Cells.Item(36, 5).Value = "  -0.75%  "

$ws.Cells.Item(37, 5).Value = "  -1.26%  "

$ws.Cells.Item(38, 5).Value = "  -8.27%  "

$ws.Cells.Item(39, 5).Value = "  +0.02%  "

$ws.Cells.Item(40, 5).Value = "  -1.28%  "

$c = $ws.Cells.Item(41, 4)
$c.Value = "'3.54"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.85%  "

$c = $ws.Cells.Item(42, 4)
$c.Value = "'0.100"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +4.43%  "

$c = $ws.Cells.Item(43, 4)
$c.Value = "'0.994"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.14%  "

$c = $ws.Cells.Item(44, 4)
$c.Value = "'276.99"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +3.54%  "

$c = $ws.Cells.Item(45, 4)
$c.Value = "'0.600"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -1.71%  "

$ws.Cells.Item(46, 5).Value = "  -4.40%  "

$c = $ws.Cells.Item(47, 4)
$c.Value = "'10.26"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.08%  "

$c = $ws.Cells.Item(48, 4)
$c.Value = "'0.0230"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -0.34%  "

$c = $ws.Cells.Item(49, 4)
$c.Value = "'4.65"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -4.75%  "

$c = $ws.Cells.Item(50, 4)
$c.Value = "'17.92"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.01%  "

$ws.Cells.Item(51, 4).Value = "1.888.39"
$ws.Cells.Item(51, 5).Value = "  +1.94%  "
